# MYCE version employed for SDEWES
# Collapse the "Renewable Sources Data" sheet from the 3-upgrade layout down
# to a single-upgrade layout, and fill in the computed results for the
# upgrade-1 figures (Units / Total Nominal Capacity / Investment / Yearly O&M
# Cost), matching a single-year optimisation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "upgrade 2" / "upgrade 3" rows (Units, Total Nominal Capacity,
# Investment, Yearly O&M Cost at upgrade 2 & 3), working from the bottom up
# so earlier row numbers stay valid while deleting.
$ws.Rows("16:17").Delete()   # Yearly O&M Cost at upgrade 2 / 3
$ws.Rows("13:14").Delete()   # Investment at upgrade 2 / 3
$ws.Rows("10:11").Delete()   # Total Nominal Capacity at upgrade 2 / 3
$ws.Rows("7:8").Delete()     # Units at upgrade 2 / 3

# Fill in the computed "Source 1" results for the remaining upgrade-1 rows.
$ws.Range("B6").Value = 45.897096546460141    # Units at upgrade 1
$ws.Range("B7").Value = 45897.09654646014     # Total Nominal Capacity at upgrade 1
$ws.Range("B8").Value = 215427.20206011989    # Investment at upgrade 1
$ws.Range("B9").Value = 4308.5440412023991    # Yearly O&M Cost at upgrade 1

# Column B now needs to be wider to fit the larger numbers, while column C
# (all zeros) can be a touch narrower - mirrors Excel's own best-fit resize.
$ws.Columns("B").ColumnWidth = 11
$ws.Columns("C").ColumnWidth = 7.1667
